$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp in the report header (A1, merged A1:D1).
$ws.Range("A1").Value = "Protractor results for: 2016-11-03 17:21:17`n"
$ws.Rows.Item(1).AutoFit()

# New rows describing the run configuration, suite, and individual specs
# (this is what lets screenshots/results be correlated with a suite+spec).
$ws.Range("A2").Value = "AppDir:./"

$ws.Range("A3").Value = "Suite:"
$ws.Range("B3").Value = "QuickStart E2E Tests"
$ws.Range("C3").Value = "passed"

$ws.Range("B4").Value = "passed"
$ws.Range("C4").Value = "should open index page"
$ws.Range("C4").WrapText = $true

$ws.Range("B5").Value = "passed"
$ws.Range("C5").Value = "should display app.component heading"
$ws.Range("C5").WrapText = $true
